$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week2")

$ws.Range("C10").Value = 0.013888888888888888
$ws.Range("C11").Value = 0.013888888888888888
$ws.Range("C15").Value = 0.010416666666666666

$ws.Range("C12").Select()
